# Daily attendance processing - 2026-01-24 16:39:17
# For every "Recorded By" cell (column G) that lists "System" and the
# instructor email together, swap the order so the email comes first,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$firstRow = $used.Row
$rowCount = $used.Rows.Count
$changed = 0

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Output "Updated $changed 'Recorded By' cell(s) in column G"
